# Updated Instruction Set Again to include Parity Check
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: fill in C27 (SUB) and E27 (SHR) that used to live in row 28
$ws.Range("C27").Value = "SUB"
$ws.Range("E27").Value = "SHR"

# Row 28: replace the old SUB entry with the new PAR (Parity Check) instruction
# and clear out the now-unused D28 cell (its value moved up to E27)
$ws.Range("B28").Value = "PAR"
$ws.Range("D28").ClearContents()

# Update the selected cell to reflect where the edit was made
$ws.Range("E27").Select()
